$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: GSW @ HOU, away, 2025-04-30 (GM5 of the series)
$row10 = @(8, "GSW", "HOU", "away", "2025-04-30", "240:00", 43, 103, 0.417, 15, 44, 0.341, 15, 22, 0.6820000000000001, 25, 24, 49, 25, 10, 5, 14, 27, 116, -15, 24, 25, 31, 36, "L")

# Row 11: HOU vs GSW, home, 2025-04-30
$row11 = @(9, "HOU", "GSW", "home", "2025-04-30", "240:00", 43, 78, 0.551, 13, 30, 0.433, 32, 38, 0.842, 8, 31, 39, 23, 9, 7, 12, 22, 131, 15, 40, 36, 31, 24, "W")

for ($i = 0; $i -lt $row10.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(10, $col).Value = $row10[$i]
    $ws.Cells.Item(11, $col).Value = $row11[$i]
}

# Excel auto-detects the DATE column values as dates; force them back to plain text
# like the rest of the sheet (matching columns B/F etc. stored as strings).
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "2025-04-30"
$ws.Cells.Item(10, 5).ClearFormats()

$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "2025-04-30"
$ws.Cells.Item(11, 5).ClearFormats()

# Column A carries the bold/bordered header-like style used throughout the table;
# copy it from the row above so the new rows match the existing formatting.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)

# Re-apply values for A10/A11 and the date cells since PasteSpecial(formats) and
# ClearFormats only touch formatting, not the already-set values.
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
